$d = $word.ActiveDocument

# 1. RF01 paragraph 1: "tiene que haber un logo de la UAH" -> "tiene estar el logo del blog, el cual fue suministrado previamente"
$d.Content.Find.Execute("En el <header> del popup tiene que haber un logo de la UAH.", $true, $false, $false, $false, $false, $true, 1, $false, "En el <header> del popup tiene estar el logo del blog, el cual fue suministrado previamente.", 2)

# 2. RF01 paragraph 2: rewritten sentence about clicking the logo
$d.Content.Find.Execute(" En el cual si nosotros pinchamos en él nos tiene que redirigir a la página del blog.", $true, $false, $false, $false, $false, $true, 1, $false, " En dicho logo, si nosotros pinchamos, nos tiene que redirigir a la página web del blog.", 2)

# 3. RF02: "Debajo del logo de la UAH ... en el blog." -> "...en la página web del blog."
$d.Content.Find.Execute("Debajo del logo de la UAH tiene que haber un buscador, para que busque otras entradas en el blog.", $true, $false, $false, $false, $false, $true, 1, $false, "Debajo del logo tiene que haber un buscador, para que busque otras entradas en la página web del blog.", 2)

# 4. 2.2 heading text stays the same, but need to insert _GoBack bookmark into RF03 paragraph
#    and remove the old one from the heading. First, move cursor to right after "blog" (before the final period)
#    in the RF03 paragraph, insert bookmark there.
$r = $d.Content
$r.Find.Execute("mostrar los últimos 5 post del blog", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)

# 5. Remove the old _GoBack bookmark in the heading (if it still exists after moving)
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    if ($bm.Range.Start -lt 2000) {
        $bm.Delete()
    }
}

# 6. Heading "2.2 Otros requisitos" -- merge runs / drop lastRenderedPageBreak
$d.Content.Find.Execute("2.2 Otros requisitos", $true, $false, $false, $false, $false, $true, 1, $false, "2.2 Otros requisitos", 2)

# 7. RNF04 "tamaño" paragraph -- merge runs (no text change)
$d.Content.Find.Execute("El tamaño de letra a utilizar en el popup tiene que ser de ", $true, $false, $false, $false, $false, $true, 1, $false, "El tamaño de letra a utilizar en el popup tiene que ser de ", 2)

# 8. "El color de letra..." paragraph -- merge runs (no text change)
$d.Content.Find.Execute("El color de letra a utilizar en el popup tiene que ser negro.", $true, $false, $false, $false, $false, $true, 1, $false, "El color de letra a utilizar en el popup tiene que ser negro.", 2)

# 9. "La anchura del popup tiene que ser de " -- merge runs (no text change)
$d.Content.Find.Execute("La anchura del popup tiene que ser de ", $true, $false, $false, $false, $false, $true, 1, $false, "La anchura del popup tiene que ser de ", 2)

# 10. "Al pasar el cursor por alguno de los iconos de las redes sociales, se tiene que mostrar un mensaje." -- merge runs (no text change)
$d.Content.Find.Execute("Al pasar el cursor por alguno de los iconos de las redes sociales, se tiene que mostrar un mensaje.", $true, $false, $false, $false, $false, $true, 1, $false, "Al pasar el cursor por alguno de los iconos de las redes sociales, se tiene que mostrar un mensaje.", 2)

# 11. RF04 footer paragraph -- merge runs / drop proofErr spell-check tags (no text change)
$d.Content.Find.Execute(" del popup, en el <footer>, en el cual si nosotros en cualquiera de los iconos nos tiene que redirigir a la red social seleccionada.", $true, $false, $false, $false, $false, $true, 1, $false, " del popup, en el <footer>, en el cual si nosotros en cualquiera de los iconos nos tiene que redirigir a la red social seleccionada.", 2)
